$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.497
$ws.Cells.Item(2, 3).Value = "BALLENOIL"
$ws.Cells.Item(2, 4).Value = "AVENIDA FUENLABRADA, 6"
$ws.Cells.Item(2, 5).Value = "HUMANES DE MADRID"

# Row 3
$ws.Cells.Item(3, 2).Value = 1.514
$ws.Cells.Item(3, 3).Value = "SHELL"
$ws.Cells.Item(3, 4).Value = "AVENIDA FUENLABRADA, 110"
$ws.Cells.Item(3, 5).Value = "HUMANES DE MADRID"

# Row 4
$ws.Cells.Item(4, 2).Value = 1.599
$ws.Cells.Item(4, 3).Value = "BALLENOIL"
$ws.Cells.Item(4, 4).Value = "CALLE CABO RUFINO LAZARO, 7"
$ws.Cells.Item(4, 5).Value = "ROZAS DE MADRID (LAS)"

# Row 5
$ws.Cells.Item(5, 2).Value = 1.649
$ws.Cells.Item(5, 3).Value = "CARREFOUR"
$ws.Cells.Item(5, 4).Value = "CARRETERA MADRID-LA CORUÑA KM. 22"
$ws.Cells.Item(5, 5).Value = "ROZAS DE MADRID (LAS)"

# Row 6
$ws.Cells.Item(6, 2).Value = 1.669
$ws.Cells.Item(6, 3).Value = "GALP"
$ws.Cells.Item(6, 4).Value = "CTRA. N-VI km 21,700"
$ws.Cells.Item(6, 5).Value = "ROZAS DE MADRID (LAS)"

# Row 7
$ws.Cells.Item(7, 2).Value = 1.688
$ws.Cells.Item(7, 3).Value = "REPSOL"
$ws.Cells.Item(7, 4).Value = "CALLE COPENHAGUES/N, S/N"
$ws.Cells.Item(7, 5).Value = "ROZAS DE MADRID (LAS)"

# Row 8
$ws.Cells.Item(8, 2).Value = 1.688
$ws.Cells.Item(8, 3).Value = "REPSOL"
$ws.Cells.Item(8, 4).Value = "A-6 km 25,5"
$ws.Cells.Item(8, 5).Value = "ROZAS DE MADRID (LAS)"

# Row 9
$ws.Cells.Item(9, 2).Value = 1.688
$ws.Cells.Item(9, 3).Value = "REPSOL"
$ws.Cells.Item(9, 4).Value = "CTRA. M-505 km 5,5"
$ws.Cells.Item(9, 5).Value = "ROZAS DE MADRID (LAS)"

# Row 10
$ws.Cells.Item(10, 2).Value = 1.688
$ws.Cells.Item(10, 3).Value = "REPSOL"
$ws.Cells.Item(10, 4).Value = "CARRETERA M-505 km 5.5"
$ws.Cells.Item(10, 5).Value = "ROZAS DE MADRID (LAS)"

# Row 11
$ws.Cells.Item(11, 2).Value = 1.689
$ws.Cells.Item(11, 3).Value = "CEPSA"
$ws.Cells.Item(11, 4).Value = "CARRETERA M-405 KM. 6"
$ws.Cells.Item(11, 5).Value = "HUMANES DE MADRID"

# Row 12
$ws.Cells.Item(12, 2).Value = 1.699
$ws.Cells.Item(12, 3).Value = "REPSOL HUMANES"
$ws.Cells.Item(12, 4).Value = "AVENIDA LA INDUSTRIA, S/N"
$ws.Cells.Item(12, 5).Value = "HUMANES DE MADRID"

# Row 13
$ws.Cells.Item(13, 2).Value = 1.699
$ws.Cells.Item(13, 3).Value = "CEPSA"
$ws.Cells.Item(13, 4).Value = "CARRETERA M-405 KM. 5,6"
$ws.Cells.Item(13, 5).Value = "HUMANES DE MADRID"

# Row 14
$ws.Cells.Item(14, 2).Value = 1.699
$ws.Cells.Item(14, 3).Value = "REPSOL"
$ws.Cells.Item(14, 4).Value = "CL MADRID, 52"
$ws.Cells.Item(14, 5).Value = "HUMANES DE MADRID"

# Row 15
$ws.Cells.Item(15, 2).Value = 1.709
$ws.Cells.Item(15, 3).Value = "BP VALDONAIRE"
$ws.Cells.Item(15, 4).Value = "CARRETERA AVD.DE LA INDUSTRIA KM. 15"
$ws.Cells.Item(15, 5).Value = "HUMANES DE MADRID"

# Row 16
$ws.Cells.Item(16, 2).Value = 1.709
$ws.Cells.Item(16, 3).Value = "BP HUMANES - EL MOLINO"
$ws.Cells.Item(16, 4).Value = "AVENIDA DE LAS FLORES, 2"
$ws.Cells.Item(16, 5).Value = "HUMANES DE MADRID"

# Row 17
$ws.Cells.Item(17, 2).Value = 1.709
$ws.Cells.Item(17, 3).Value = "REPSOL"
$ws.Cells.Item(17, 4).Value = "CARRETERA AVENIDA  DE LA INDUSTRIA , 46 KM. 1,1"
$ws.Cells.Item(17, 5).Value = "HUMANES DE MADRID"

# Row 18
$ws.Cells.Item(18, 2).Value = 1.709
$ws.Cells.Item(18, 3).Value = "REPSOL"
$ws.Cells.Item(18, 4).Value = "CR A-6, 20,3"
$ws.Cells.Item(18, 5).Value = "ROZAS DE MADRID (LAS)"

# Row 19
$ws.Cells.Item(19, 2).Value = 1.709
$ws.Cells.Item(19, 3).Value = "BP LAS ROZAS"
$ws.Cells.Item(19, 4).Value = "CL LAS CRUCES  S/N"
$ws.Cells.Item(19, 5).Value = "ROZAS DE MADRID (LAS)"

# Row 20
$ws.Cells.Item(20, 2).Value = 1.861
$ws.Cells.Item(20, 3).Value = "COSTCO"
$ws.Cells.Item(20, 4).Value = "CALLE INNOVACIÓN, 19"
$ws.Cells.Item(20, 5).Value = "ROZAS DE MADRID (LAS)"

# Row 21
$ws.Cells.Item(21, 3).Value = "T9"
$ws.Cells.Item(21, 4).Value = "CALLE TENERIFE (DE), 2"
$ws.Cells.Item(21, 5).Value = "HUMANES DE MADRID"

Write-Host "applied price/rotulo/direccion/localidad updates"